$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 661 entirely (the post about "彼は完璧ではない...") — all rows
# below it shift up by one automatically.
$ws.Rows(661).Delete()
